$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting existing rows 68:86 down to 69:87
$ws.Rows.Item(68).Insert()

# Populate the new row 68 - it mirrors the surrounding "Ají" records for this market,
# differing only in date, variety, volume, prices, origin and price/kg.
$ws.Range("A68").Value = 7
$ws.Range("B68").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C68").Value = "Ñuble"
$ws.Range("D68").Value = 44754
$ws.Range("E68").Value = 16
$ws.Range("F68").Value = 100112021
$ws.Range("G68").Value = "Ají"
$ws.Range("H68").Value = "Inferno"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 80
$ws.Range("K68").Value = 18000
$ws.Range("L68").Value = 19000
$ws.Range("M68").Value = 18500
$ws.Range("N68").Value = "$/caja 15 kilos"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 1233
$ws.Range("Q68").Value = 15
$ws.Range("R68").Value = "Hortaliza"
